$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value that "looks like a date" (e.g. "2026-02-25")
# without letting Excel auto-convert it to a date serial. A leading
# apostrophe forces text entry (quotePrefix); re-applying the "Normal"
# style afterwards clears the quote-prefix style bit again so the cell
# ends up with no explicit style, just like the surrounding data cells.
function Set-TextValue($rangeAddr, [string]$value) {
    $rng = $ws.Range($rangeAddr)
    $rng.Formula = "'" + $value
    $rng.Style = "Normal"
}

# --- Row 2 ---
$ws.Range("A2").Value = "Software Engineer II - Content Platform Engineering"
$ws.Range("B2").Value = "nan"
$ws.Range("C2").Value = "Bristol, CT, US USA"
$ws.Range("D2").Value = 18.9
$ws.Range("E2").Value = "Data Scientist, Generative AI, RAG, Prompt Engineering, TensorFlow, PyTorch, S3, EC2, Docker, Kubernetes"
Set-TextValue "F2" "2026-02-25"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=63b62b8246a9f659"

# --- Row 3 ---
$ws.Range("A3").Value = "Machine Learning Engineer"
$ws.Range("B3").Value = "CVS Health"
$ws.Range("D3").Value = 16.7
$ws.Range("E3").Value = "Data Scientist, Machine Learning Engineer, RAG, TensorFlow, PyTorch, BigQuery, Snowflake, BigQuery, Kafka, Python"
Set-TextValue "F3" "2026-02-25"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=f836fc84892ec33f"

# --- Row 4 ---
$ws.Range("A4").Value = "Data Scientist - Pricing and Promotions Optimization"
$ws.Range("B4").Value = "CVS Health"
$ws.Range("D4").Value = 13.3
$ws.Range("E4").Value = "Data Scientist, AWS SageMaker, GCP Vertex AI, Git, Databricks, PySpark, Hadoop, Python, SQL, R"
Set-TextValue "F4" "2026-02-25"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=f2dd0084730de1d4"

# --- Row 5 ---
$ws.Range("A5").Value = "Data Scientist - Pricing and Promotions Optimization"
$ws.Range("B5").Value = "CVS Health"
$ws.Range("C5").Value = "Wellesley, MA, US USA"
$ws.Range("D5").Value = 13.3
$ws.Range("E5").Value = "Data Scientist, AWS SageMaker, GCP Vertex AI, Git, Databricks, PySpark, Hadoop, Python, SQL, R"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=7ac1b66db18d2c42"

# --- Row 6 ---
$ws.Range("A6").Value = "Software Engineer - DevOps"
$ws.Range("B6").Value = "nan"
$ws.Range("C6").Value = "New York, NY, US USA"
$ws.Range("D6").Value = 13.3
$ws.Range("E6").Value = "Data Scientist, Copilot, Docker, CI/CD, Jenkins, GitHub Actions, Terraform, Git, Python, R"
Set-TextValue "F6" "2026-02-25"
$ws.Range("G6").Value = "https://www.indeed.com/viewjob?jk=c10fe2dd154808d6"

# --- Row 7 ---
$ws.Range("A7").Value = "Data Scientist"
$ws.Range("B7").Value = "Rice University"
$ws.Range("C7").Value = "Houston, TX, US USA"
$ws.Range("E7").Value = "Data Scientist, Git, Hadoop, Tableau, Power BI, Python, SQL, R, Scala, Optimization"
Set-TextValue "F7" "2026-02-25"
$ws.Range("G7").Value = "https://www.indeed.com/viewjob?jk=29c39de545706a17"

# --- Row 8 ---
$ws.Range("A8").Value = "AI Solutions Engineer"
$ws.Range("B8").Value = "Conagra Brands"
$ws.Range("C8").Value = "Chicago, IL, US USA"
$ws.Range("E8").Value = "Generative AI, LangChain, RAG, FAISS, Pinecone, Python, SQL, R, Java, Optimization"
Set-TextValue "F8" "2026-02-25"
$ws.Range("G8").Value = "https://www.indeed.com/viewjob?jk=1b181f522ed1f1b9"

# --- Row 9 ---
$ws.Range("A9").Value = "AI Engineering Intern (Generative & Agentic AI)"
$ws.Range("B9").Value = "Boston Scientific"
$ws.Range("C9").Value = "Marlboro, MA, US USA"
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = "AI Engineer, Generative AI, LangChain, RAG, Prompt Engineering, Python, R, Scala, Optimization"
Set-TextValue "F9" "2026-02-25"
$ws.Range("G9").Value = "https://www.indeed.com/viewjob?jk=7682f4c2d9fcb8d0"

# --- Remove old row 10 entirely (shrinks used range to A1:G9) ---
$ws.Rows.Item(10).Delete()
